# Insert a new bullet right after the "Created the first view ... to show
# the list of Items on the application." paragraph, per the commit:
#   "Created the ViewModels folder and added the first view model for
#    Items and passed Item content to the List.chtml template using the
#    ItemController."
# The new paragraph reuses the same ListParagraph/numbering formatting as
# its sibling bullets, and its text is split across four runs (all with
# identical <w:lang w:val="en-US"/> run formatting) to mirror the target
# OOXML exactly.

$d = $word.ActiveDocument

# Locate the anchor paragraph (the last bullet about the "first view").
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*to show the list of Items on the application.*") {
        $anchorPara = $candidate
        break
    }
}
if ($null -eq $anchorPara) {
    throw "Could not locate anchor paragraph ending in '...application.'"
}

# Create a brand new (empty) paragraph right after it.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newRange = $newPara.Range
$newRange.Collapse(1)   # wdCollapseStart

# Inject the bullet's full OOXML (pPr + four w:r runs) via InsertXML so the
# run boundaries match the target exactly instead of being coalesced into a
# single run by plain text insertion.
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Created the ViewModels folder and added the first view model for Items and passed</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Item</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> content to the List.chtml template </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>using the ItemController.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part>
</pkg:package>
'@

$newRange.InsertXML($xml)
